$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the existing header style (bold, border, centered) onto the three new
# header cells before writing their text, so AD1:AF1 look like the rest of row 1.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122) # xlPasteFormats

# New headers for team win/loss/tie record
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Every player row (2-49) gets the same team record values
for ($r = 2; $r -le 49; $r++) {
    $ws.Cells.Item($r, 30).Value = 75  # AD: Wins
    $ws.Cells.Item($r, 31).Value = 87  # AE: Losses
    $ws.Cells.Item($r, 32).Value = 0   # AF: Ties
}
